# Submit the "Is a solar future inevitable version"
# Model version used for the Nature Communications paper
# "Is a solar future inevitable".

$wb = $excel.ActiveWorkbook

# --- FTT-P sheet: just a cursor/selection move ---
$wsP = $wb.Worksheets.Item("FTT-P")
$wsP.Range("H13").Select() | Out-Null

# --- FTT-Fr sheet: mark rows as active (Column B 0 -> 1) and drop the
#     now-retired last data row (row 27, "ZCEZ" / T-Scaling placeholder) ---
$wsFr = $wb.Worksheets.Item("FTT-Fr")

$rowsToFlip = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,18,19,20,21,22,23,24,26)
foreach ($r in $rowsToFlip) {
    $wsFr.Cells.Item($r, 2).Value = 1
}

# Remove the trailing row (row 27) entirely - this also drops its four
# now-unused shared strings (ZCEZ / None / Costs (Column 64) / FTT-Fr T-Scaling)
$wsFr.Rows.Item(27).Delete()

# Move the selection on FTT-Fr away from the old C31 cursor
$wsFr.Range("A7").Select() | Out-Null

# --- Time_Horizons sheet becomes the active tab/sheet ---
$wsTH = $wb.Worksheets.Item("Time_Horizons")
$wsTH.Activate()
